$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.434.00"
$ws.Range("E2").Value = "  +0.94%  "

# Row 3
$ws.Range("D3").Value = "1.796.78"
$ws.Range("E3").Value = "  +0.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.554"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.40%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.53"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.296"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.31%  "

# Row 10
$ws.Range("E10").Value = "  +0.40%  "

# Row 11
$ws.Range("E11").Value = "  +0.62%  "

# Row 12
$ws.Range("D12").Value = "2.055.68"
$ws.Range("E12").Value = "  +0.56%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.11%  "

# Row 14
$ws.Range("D14").Value = "1.790.92"
$ws.Range("E14").Value = "  +0.48%  "

# Row 16
$ws.Range("D16").Value = "34.388.09"
$ws.Range("E16").Value = "  +0.89%  "

# Row 17
$ws.Range("E17").Value = "  +1.08%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.61%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0802"
$ws.Range("E19").Value = "  +3.07%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "246.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.44%  "

# Row 22
$ws.Range("E22").Value = "  +0.12%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.80%  "

# Row 24
$ws.Range("E24").Value = "  +0.53%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "

# Row 26
$ws.Range("E26").Value = "  +0.65%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.66%  "

# Row 28
$ws.Range("E28").Value = "  +1.94%  "

# Row 29
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.44%  "

# Row 31
$ws.Range("E31").Value = "  +0.91%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.44%  "

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.78%  "

# Row 34
$ws.Range("E34").Value = "  +1.33%  "

# Row 35
$ws.Range("D35").Value = "1.444.52"
$ws.Range("E35").Value = "  -0.42%  "

# Row 36
$ws.Range("E36").Value = "  +9.81%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.668"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.46%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.86%  "

# Row 39
$ws.Range("E39").Value = "  -0.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.77%  "

# Row 41
$ws.Range("E41").Value = "  +1.34%  "

# Row 42
$ws.Range("E42").Value = "  +3.04%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.935"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.90%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.61%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0525"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.53%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "

# Row 47
$ws.Range("E47").Value = "  +0.07%  "

# Row 48
$ws.Range("D48").Value = "1.951.61"
$ws.Range("E48").Value = "  +0.31%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.38%  "

# Row 50
$ws.Range("E50").Value = "  +0.12%  "

# Row 51
$ws.Range("D51").Value = "0.0₆0130"
$ws.Range("E51").Value = "  -4.59%  "
